# Apply the scheduled market-data refresh to the Leve profit sheets.
# Each sheet (one per crafting class) gets updated currentAveragePrice /
# Leve price / profit columns (H-N) for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 180.25
$ws.Range("I5").Value = 110.4
$ws.Range("J5").Value = 296.66666
$ws.Range("K5").Value = 110.4
$ws.Range("L5").Value = 296.66666
$ws.Range("M5").Value = 4.599999999999994
$ws.Range("N5").Value = -526.66666
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 5
$ws.Range("K6").Value = 15
$ws.Range("M6").Value = 97
$ws.Range("H9").Value = 1413.7778
$ws.Range("I9").Value = 1718
$ws.Range("K9").Value = 1718
$ws.Range("M9").Value = -1549
$ws.Range("H33").Value = 165.88889
$ws.Range("I33").Value = 165.88889
$ws.Range("K33").Value = 165.88889
$ws.Range("M33").Value = 63.11111
$ws.Range("H40").Value = 2235.3635
$ws.Range("I40").Value = 2041.7142
$ws.Range("K40").Value = 2041.7142
$ws.Range("M40").Value = -1866.7142
$ws.Range("H99").Value = 261.75
$ws.Range("J99").Value = 215
$ws.Range("L99").Value = 645
$ws.Range("N99").Value = -3641
$ws.Range("H116").Value = 5224.75
$ws.Range("I116").Value = 3999.6667
$ws.Range("K116").Value = 3999.6667
$ws.Range("M116").Value = -557.6667000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2704.3635
$ws.Range("I61").Value = 2704.3635
$ws.Range("K61").Value = 2704.3635
$ws.Range("M61").Value = -2492.3635
$ws.Range("H122").Value = 1920.5555
$ws.Range("I122").Value = 1920.5555
$ws.Range("K122").Value = 5761.666499999999
$ws.Range("M122").Value = -3311.666499999999
$ws.Range("H132").Value = 2584.3
$ws.Range("I132").Value = 2567.6843
$ws.Range("K132").Value = 7703.0529
$ws.Range("M132").Value = -5173.0529
$ws.Range("H136").Value = 2704.3635
$ws.Range("I136").Value = 2704.3635
$ws.Range("K136").Value = 8113.0905
$ws.Range("M136").Value = -5563.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 910
$ws.Range("J64").Value = 889.4
$ws.Range("L64").Value = 889.4
$ws.Range("N64").Value = -1339.4
$ws.Range("H67").Value = 910
$ws.Range("J67").Value = 889.4
$ws.Range("L67").Value = 889.4
$ws.Range("N67").Value = -2449.4
$ws.Range("H96").Value = 4625
$ws.Range("I96").Value = 4625
$ws.Range("K96").Value = 4625
$ws.Range("M96").Value = -1879
$ws.Range("H99").Value = 1999.8334
$ws.Range("I99").Value = 1999.8334
$ws.Range("K99").Value = 1999.8334
$ws.Range("M99").Value = -501.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1449.5
$ws.Range("I31").Value = 1426.75
$ws.Range("J31").Value = 1495
$ws.Range("K31").Value = 1426.75
$ws.Range("L31").Value = 1495
$ws.Range("M31").Value = -1131.75
$ws.Range("N31").Value = -2085
$ws.Range("H34").Value = 1449.5
$ws.Range("I34").Value = 1426.75
$ws.Range("J34").Value = 1495
$ws.Range("K34").Value = 1426.75
$ws.Range("L34").Value = 1495
$ws.Range("M34").Value = -1224.75
$ws.Range("N34").Value = -1899
$ws.Range("H62").Value = 1799
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1799
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 1799
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3047
$ws.Range("H65").Value = 1799
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1799
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 8995
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -15235
$ws.Range("H99").Value = 3370
$ws.Range("I99").Value = 3365.7144
$ws.Range("J99").Value = 3400
$ws.Range("K99").Value = 3365.7144
$ws.Range("L99").Value = 3400
$ws.Range("M99").Value = -1867.7144
$ws.Range("N99").Value = -6396
$ws.Range("H100").Value = 175390
$ws.Range("J100").Value = 175390
$ws.Range("L100").Value = 175390
$ws.Range("N100").Value = -177554
$ws.Range("H126").Value = 3370
$ws.Range("I126").Value = 3365.7144
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 10097.1432
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -7627.143199999999
$ws.Range("N126").Value = -15140
$ws.Range("H134").Value = 2015.4615
$ws.Range("I134").Value = 1991
$ws.Range("K134").Value = 5973
$ws.Range("M134").Value = -3438

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 128.4
$ws.Range("I17").Value = 127.666664
$ws.Range("K17").Value = 382.999992
$ws.Range("M17").Value = -213.999992
$ws.Range("H51").Value = 1423.5
$ws.Range("I51").Value = 483.33334
$ws.Range("K51").Value = 1450.00002
$ws.Range("M51").Value = -990.0000199999999
$ws.Range("H113").Value = 1557.1818
$ws.Range("J113").Value = 1536.25
$ws.Range("L113").Value = 4608.75
$ws.Range("N113").Value = -8948.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3649.375
$ws.Range("J80").Value = 3516.1667
$ws.Range("L80").Value = 3516.1667
$ws.Range("N80").Value = -5512.1667
$ws.Range("H83").Value = 3649.375
$ws.Range("J83").Value = 3516.1667
$ws.Range("L83").Value = 17580.8335
$ws.Range("N83").Value = -27564.8335
$ws.Range("H122").Value = 2489.2222
$ws.Range("I122").Value = 2203.6667
$ws.Range("J122").Value = 2632
$ws.Range("K122").Value = 6611.000100000001
$ws.Range("L122").Value = 7896
$ws.Range("M122").Value = -4161.000100000001
$ws.Range("N122").Value = -12796
$ws.Range("H132").Value = 1955.5555
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8882.5
$ws.Range("I7").Value = 3300
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 3300
$ws.Range("L7").Value = 9999
$ws.Range("M7").Value = -3188
$ws.Range("N7").Value = -10223
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H46").Value = 2288.4
$ws.Range("J46").Value = 3750
$ws.Range("L46").Value = 3750
$ws.Range("N46").Value = -4126
$ws.Range("H55").Value = 889.1667
$ws.Range("I55").Value = 633.6
$ws.Range("J55").Value = 2167
$ws.Range("K55").Value = 633.6
$ws.Range("L55").Value = 2167
$ws.Range("M55").Value = -460.6
$ws.Range("N55").Value = -2513
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H126").Value = 8882.5
$ws.Range("I126").Value = 3300
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 9900
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -7430
$ws.Range("N126").Value = -34937

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 24282.715
$ws.Range("I14").Value = 28327.334
$ws.Range("K14").Value = 28327.334
$ws.Range("M14").Value = -28159.334
$ws.Range("H46").Value = 44245
$ws.Range("J46").Value = 44245
$ws.Range("L46").Value = 44245
$ws.Range("N46").Value = -44707
$ws.Range("H107").Value = 683.5454999999999
$ws.Range("J107").Value = 923.5
$ws.Range("L107").Value = 2770.5
$ws.Range("N107").Value = -6610.5
$ws.Range("H126").Value = 5412.1724
$ws.Range("J126").Value = 5665
$ws.Range("L126").Value = 16995
$ws.Range("N126").Value = -21935
$ws.Range("H132").Value = 1522.72
$ws.Range("I132").Value = 1480.091
$ws.Range("K132").Value = 4440.272999999999
$ws.Range("M132").Value = -1910.272999999999
$ws.Range("H134").Value = 44245
$ws.Range("J134").Value = 44245
$ws.Range("L134").Value = 132735
$ws.Range("N134").Value = -137805
